$d = $word.ActiveDocument

# 1. Insert a new Title paragraph before the existing first paragraph
#    ("This is the first page.").
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertBefore("This is my doc, hehe`r")
$d.Paragraphs.Item(1).Style = $d.Styles.Item("Title")

# 2. Insert a new Heading 1 paragraph right after "This is the first page."
#    (now paragraph 2) and before the following blank paragraph, without
#    disturbing that blank paragraph.
$followingPara = $d.Paragraphs.Item(3)
$followingPara.Range.InsertBefore("Random heading lol 1`r")
$d.Paragraphs.Item(3).Style = $d.Styles.Item("Heading 1")
